$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a serial date value (2023-09-09 == 45178)
# that was bumped by one day (to 2023-09-10 == 45179) for every data row
# (rows 2 through 199).
for ($r = 2; $r -le 199; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45178) {
        $cell.Value2 = 45179
    }
}
